# Generate Report for handback
#
# Populates the "Latest Target File" (E) / "Latest Handback File" (F) columns
# for the two localized files in the zh-cn and de-de sheets, refreshes the
# "Latest Handback DateTime" (G) for those rows, and flips the Status (B)
# text from "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Per-language sheet definitions: source repo URLs differ only by the
# handoff-org branch SHA and the language folder, so keep them data-driven.
# ---------------------------------------------------------------------------
$sheetsInfo = @(
    @{
        Name = "zh-cn"
        MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/54dc6125d47821a4738aab7f9c31fed998ae84a4/e2e/"
        XlfUrlPrefix = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/89f58944277784a77e96dd87122ef287f5b9db5c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/"
        XlfSuffix = "zh-cn.xlf"
        HandbackDateTime = "2016-01-18 02:50:38"
    },
    @{
        Name = "de-de"
        MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/54dc6125d47821a4738aab7f9c31fed998ae84a4/e2e/"
        XlfUrlPrefix = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d778ccf1e75e2beb28b483ec1761e3b23804e8f4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/"
        XlfSuffix = "de-de.xlf"
        HandbackDateTime = "2016-01-18 02:50:57"
    }
)

$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/54dc6125d47821a4738aab7f9c31fed998ae84a4/.localization-config"

$row2Md = "24c6999a-c474-456a-a3c4-d2ff4879f823.md"
$row2XlfBase = "24c6999a-c474-456a-a3c4-d2ff4879f823.1e20969ee52e7e1786bd75d2aa13112624e3df54"
$row3Md = "879e6e22-dc92-4e1c-823e-b54891410cd1.md"
$row3XlfBase = "879e6e22-dc92-4e1c-823e-b54891410cd1.4d879a31728026479d9ccb62366e9e37147dc90d"
$configDisplay = ".localization-config"

$statusHandedBack = "Handed back: in sync with en-US"

foreach ($info in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    $row2Xlf = "$row2XlfBase.$($info.XlfSuffix)"
    $row3Xlf = "$row3XlfBase.$($info.XlfSuffix)"

    # -- Status column: report is now handed back, in sync with en-US -------
    $ws.Range("B2").Value = $statusHandedBack
    $ws.Range("B3").Value = $statusHandedBack

    # -- New "Latest Target File" (E) / "Latest Handback File" (F) values ---
    $ws.Range("E2").Value = $row2Md
    $ws.Range("F2").Value = $row2Xlf
    $ws.Range("E3").Value = $row3Md
    $ws.Range("F3").Value = $row3Xlf

    # -- Latest Handback DateTime (G) now reflects the real handback time ---
    $ws.Range("G2").Value = $info.HandbackDateTime
    $ws.Range("G3").Value = $info.HandbackDateTime

    # -- Rebuild hyperlinks in worksheet order so relationship ids line up --
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), ($info.MdUrl + $row2Md), "", "", $row2Md)
    $ws.Hyperlinks.Add($ws.Range("C2"), ($info.XlfUrlPrefix + $row2Xlf), "", "", $row2Xlf)
    $ws.Hyperlinks.Add($ws.Range("E2"), ($info.MdUrl + $row2Md), "", "", $row2Md)
    $ws.Hyperlinks.Add($ws.Range("F2"), ($info.XlfUrlPrefix + $row2Xlf), "", "", $row2Xlf)

    $ws.Hyperlinks.Add($ws.Range("A3"), ($info.MdUrl + $row3Md), "", "", $row3Md)
    $ws.Hyperlinks.Add($ws.Range("C3"), ($info.XlfUrlPrefix + $row3Xlf), "", "", $row3Xlf)
    $ws.Hyperlinks.Add($ws.Range("E3"), ($info.MdUrl + $row3Md), "", "", $row3Md)
    $ws.Hyperlinks.Add($ws.Range("F3"), ($info.XlfUrlPrefix + $row3Xlf), "", "", $row3Xlf)

    $ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", $configDisplay)

    # -- Give the two new columns the same "HyperLink" look as the rest -----
    $ws.Range("E2").Style = "HyperLink"
    $ws.Range("F2").Style = "HyperLink"
    $ws.Range("E3").Style = "HyperLink"
    $ws.Range("F3").Style = "HyperLink"
}

Write-Host "Handback report generated for zh-cn and de-de sheets."
